# Update daily COVID-19 Valais figures for rows 640-660 (Feuil1).
# Columns: A=Date, B=Cumul cas positifs (formula), C=Nb nouveaux cas positifs,
#          D=Nb nouvelles admissions (unused), E=Patients SI, F=Patients intubés,
#          G=Patients hospitalisés hors SI, H=Total hospitalisations (formula),
#          I=Nb nouvelles sorties, J=Cumul décès (formula), K=Nb nouveaux décès (formula),
#          L=Nb nouveaux décès hôpital, M=Nb nouveaux décès extra-hospitaliers.
# B, H, J, K are shared formulas that recalc automatically once their inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L and M are formatted as Text (@). Assigning a numeric .Value directly
# to a Text-formatted cell stores it as a text string instead of a number, so for
# those two columns toggle the format to General for the write and restore it
# right after (matches how the source data is actually stored: real numbers).
function Set-NumericValue($ref, $val) {
    $rng = $ws.Range($ref)
    $fmt = $rng.NumberFormat
    $rng.NumberFormat = "General"
    $rng.Value = $val
    $rng.NumberFormat = $fmt
}

$ws.Range("C640").Value = 268

$ws.Range("C641").Value = 153

$ws.Range("E643").Value = 11
$ws.Range("F643").Value = 3

$ws.Range("E644").Value = 11
$ws.Range("F644").Value = 6

$ws.Range("E645").Value = 10
$ws.Range("F645").Value = 6

$ws.Range("C646").Value = 342
$ws.Range("E646").Value = 10
$ws.Range("F646").Value = 8

$ws.Range("C647").Value = 354
$ws.Range("E647").Value = 10
$ws.Range("F647").Value = 9

$ws.Range("C648").Value = 140
$ws.Range("E648").Value = 10
$ws.Range("F648").Value = 8

$ws.Range("E649").Value = 10
$ws.Range("F649").Value = 7
Set-NumericValue "L649" 2

$ws.Range("C650").Value = 473
$ws.Range("E650").Value = 10
$ws.Range("F650").Value = 8

$ws.Range("E651").Value = 13
$ws.Range("F651").Value = 9

$ws.Range("C652").Value = 168
$ws.Range("E652").Value = 12
$ws.Range("F652").Value = 9

$ws.Range("C653").Value = 446
$ws.Range("E653").Value = 12
$ws.Range("F653").Value = 10

$ws.Range("C654").Value = 340
$ws.Range("E654").Value = 12
$ws.Range("F654").Value = 9

$ws.Range("C655").Value = 244
$ws.Range("E655").Value = 12
$ws.Range("F655").Value = 10

$ws.Range("C656").Value = 135
$ws.Range("E656").Value = 13
$ws.Range("F656").Value = 8
Set-NumericValue "L656" 2

$ws.Range("C657").Value = 530

# Rows 658-660 were previously blank placeholders; fill in the new daily data.
$ws.Range("C658").Value = 424
$ws.Range("E658").Value = 11
$ws.Range("F658").Value = 6
$ws.Range("G658").Value = 69
Set-NumericValue "L658" 0
Set-NumericValue "M658" 0

$ws.Range("C659").Value = 257
$ws.Range("E659").Value = 10
$ws.Range("F659").Value = 6
$ws.Range("G659").Value = 70
Set-NumericValue "L659" 0
Set-NumericValue "M659" 0

$ws.Range("C660").Value = 15
$ws.Range("E660").Value = 8
$ws.Range("F660").Value = 5
$ws.Range("G660").Value = 71
Set-NumericValue "L660" 0
Set-NumericValue "M660" 0
